# Refresh the cryptocurrency table (prices, 1h volume %, and the two
# swapped row pairs) to match the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain-looking numbers as literal text
# (e.g. "41.149.38", "0.0740", "4.60"). Assigning such strings straight
# to .Value lets Excel "helpfully" reinterpret many of them as numbers
# (losing formatting / introducing floating point noise), so every
# Price cell is forced to Text format first and restored to the default
# "Normal" style afterwards, to store the exact original literal text.
$priceCells = @(
    'D2',
    'D3',
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D15',
    'D16',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D26',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D51'
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.149.38'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '2.177.92'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '253.49'
$ws.Range('E5').Value = '  +6.44%  '
$ws.Range('D6').Value = '0.626'
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('D7').Value = '67.54'
$ws.Range('E7').Value = '  -2.88%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.577'
$ws.Range('D10').Value = '37.51'
$ws.Range('E10').Value = '  +3.70%  '
$ws.Range('D11').Value = '58.62'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '0.0931'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '7.08'
$ws.Range('E13').Value = '  +9.49%  '
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '2.501.88'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '0.866'
$ws.Range('E16').Value = '  +5.17%  '
$ws.Range('D17').Value = '14.43'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '2.189.80'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '41.098.47'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('D21').Value = '6.16'
$ws.Range('E21').Value = '  +2.63%  '
$ws.Range('D22').Value = '71.71'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '231.79'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').Value = '2.03'
$ws.Range('E24').Value = '  +2.79%  '
$ws.Range('E25').Value = '  +9.18%  '
$ws.Range('D26').Value = '11.81'
$ws.Range('E26').Value = '  +22.97%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +6.33%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '168.39'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '20.55'
$ws.Range('E31').Value = '  +2.66%  '
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0740'
$ws.Range('E33').Value = '  +6.83%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.123'
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').Value = '  +7.54%  '
$ws.Range('D36').Value = '27.67'
$ws.Range('E36').Value = '  +21.46%  '
$ws.Range('D37').Value = '4.19'
$ws.Range('E37').Value = '  +11.71%  '
$ws.Range('D38').Value = '4.60'
$ws.Range('E38').Value = '  +1.61%  '
$ws.Range('D39').Value = '0.0298'
$ws.Range('E39').Value = '  +13.67%  '
$ws.Range('D40').Value = '12.65'
$ws.Range('E40').Value = '  +25.87%  '
$ws.Range('D41').Value = '2.19'
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('D42').Value = '5.68'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').Value = '64.15'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').Value = '4.98'
$ws.Range('D45').Value = '0.200'
$ws.Range('E45').Value = '  +5.63%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.101'
$ws.Range('E46').Value = '  +3.60%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '8.58'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('E49').Value = '  +5.23%  '
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').Value = '4.29'
$ws.Range('E51').Value = '  -2.87%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

